$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 118, shifting existing rows 118:143 down to 119:144.
$ws.Rows.Item(118).Insert()

# Populate the newly inserted row 118 with the new weekly price record.
$ws.Range("A118").Value = 9
$ws.Range("B118").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C118").Value = "Metropolitana"
$ws.Range("D118").Value = 45173
$ws.Range("D118").NumberFormat = $ws.Range("D119").NumberFormat
$ws.Range("E118").Value = 13
$ws.Range("F118").Value = 100114007
$ws.Range("G118").Value = "Jengibre"
$ws.Range("H118").Value = "Sin especificar"
$ws.Range("I118").Value = "Primera"
$ws.Range("J118").Value = 520
$ws.Range("K118").Value = 17000
$ws.Range("L118").Value = 18000
$ws.Range("M118").Value = 17500
$ws.Range("N118").Value = "`$/caja 13 kilos"
$ws.Range("O118").Value = "Perú"
$ws.Range("P118").Value = 1346
$ws.Range("Q118").Value = 13
$ws.Range("R118").Value = "Hortaliza"
